$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Contest 13 (row 22, "KXI vs MI") results are now in: fill in the raw
#    points for each of the six players (E/H/K/N/Q/T). The VLOOKUP formulas
#    in D/G/J/M/P/S already exist and will recompute automatically.
# ---------------------------------------------------------------------------
$ws.Range("E22").Value2 = 20
$ws.Range("H22").Value2 = 0
$ws.Range("K22").Value2 = 40
$ws.Range("N22").Value2 = 80
$ws.Range("Q22").Value2 = 60
$ws.Range("T22").Value2 = 100

# ---------------------------------------------------------------------------
# 2) A new contest row ("contest 22", SRH vs KXI) is added right below row 31
#    (the previous last templated/blank contest row). Insert a new row at 32
#    so everything below (the truly-blank spacer row, the Team/Prize header,
#    the player-name row and the Totals row) shifts down by one, matching
#    what Excel does when you right click a row header and choose Insert.
# ---------------------------------------------------------------------------
$ws.Rows("32").Insert()

# Row 22 of the table (template row 31) already carries A/B/C values - row 31
# gets the new contest's identity: contest number 22, format 1, match name.
$ws.Range("A31").Value2 = 22
$ws.Range("B31").Value2 = 1
$ws.Range("C31").Value2 = "SRH vs KXI"

# The inserted row 32 needs to become the new "blank templated" row (same
# look as the old row 31 before it got filled in): copy the formatting from
# row 30 cell-by-cell (avoids pulling in the empty gap columns that a whole
# -row copy would introduce), then re-create the VLOOKUP formulas so they
# reference row 32.
$cols = @("A","B","C","D","E","G","H","J","K","M","N","P","Q","S","T")
foreach ($col in $cols) {
    $ws.Range($col + "30").Copy()
    $ws.Range($col + "32").PasteSpecial(-4122)
}

$ws.Range("D32").Formula = "=IF(ISERROR(VLOOKUP(RANK(E32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE)),"""",VLOOKUP(RANK(E32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE))"
$ws.Range("G32").Formula = "=IF(ISERROR(VLOOKUP(RANK(H32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE)),"""",VLOOKUP(RANK(H32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE))"
$ws.Range("J32").Formula = "=IF(ISERROR(VLOOKUP(RANK(K32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE)),"""",VLOOKUP(RANK(K32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE))"
$ws.Range("M32").Formula = "=IF(ISERROR(VLOOKUP(RANK(N32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE)),"""",VLOOKUP(RANK(N32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE))"
$ws.Range("P32").Formula = "=IF(ISERROR(VLOOKUP(RANK(Q32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE)),"""",VLOOKUP(RANK(Q32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE))"
$ws.Range("S32").Formula = "=IF(ISERROR(VLOOKUP(RANK(T32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE)),"""",VLOOKUP(RANK(T32, (`$T32,`$Q32,`$N32,`$K32,`$H32,`$E32), 0),  score, 2, FALSE))"

# ---------------------------------------------------------------------------
# 3) The season Totals row (old row 35, now row 36) needs its SUM ranges
#    extended by one row (D10:D31 -> D10:D32, etc.) to include the newly
#    inserted contest row.
# ---------------------------------------------------------------------------
$ws.Range("E36").Formula = "=SUM(D10:D32)"
$ws.Range("H36").Formula = "=SUM(G10:G32)"
$ws.Range("K36").Formula = "=SUM(J10:J32)"
$ws.Range("N36").Formula = "=SUM(M10:M32)"
$ws.Range("Q36").Formula = "=SUM(P10:P32)"
$ws.Range("T36").Formula = "=SUM(S10:S32)"

# ---------------------------------------------------------------------------
# 4) The conditional formatting (win/loss/even colouring) that used to sit on
#    the Totals row cells (E35/H35/K35/N35/Q35/T35) needs to move down to the
#    same cells on their new row (E36/H36/K36/N36/Q36/T36), keeping the exact
#    same rules/colours.
# ---------------------------------------------------------------------------
$totalsCols = @("E","H","K","N","Q","T")
foreach ($col in $totalsCols) {
    $src = $ws.Range($col + "35")
    $dstAddr = $col + "36"
    $fcCount = $src.FormatConditions.Count
    for ($i = 1; $i -le $fcCount; $i++) {
        $fc = $src.FormatConditions.Item(1)
        $fc.ModifyAppliesToRange($ws.Range($dstAddr))
    }
}

# ---------------------------------------------------------------------------
# 5) Keep the selection in sync with where Excel would have left the cursor
#    (on the grand-total cell, which is now U36 instead of U35).
# ---------------------------------------------------------------------------
$ws.Range("U36").Select()
